# Prezenta Algoritmi - add "săpt. 12" (column N) attendance for several students
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add attendance value 2 in column N for the rows below (students present in week 12)
$rows = @(4, 7, 10, 12, 13, 20, 22)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 14).Value = 2   # column N = 14
}

# Update the active cell selection on the sheet (moved from M20 to M21)
$ws.Range("M21").Select()
